$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8 - start from a copy of row 7's formatting (A:E only, row 7 has no F-less
# equivalent so copy just the columns that are populated on the new rows)
$ws.Range("A7:E7").Copy($ws.Range("A8:E8"))
$ws.Cells.Item(8, 1).Value = 5
$ws.Cells.Item(8, 3).Value = 82981004249
$ws.Cells.Item(8, 5).Value = "Bruno"

# New row 9 - same approach
$ws.Range("A7:E7").Copy($ws.Range("A9:E9"))
$ws.Cells.Item(9, 1).Value = 6
$ws.Cells.Item(9, 3).Value = 82981004249
$ws.Cells.Item(9, 5).Value = "Ronaldo"

$ws.Range("C9").Select() | Out-Null
